# "fix special label for artists"
# The EventLong (duration) values for two "Трансляция" (broadcast) rows were
# wrongly left at "00:00:00" and should be "01:00:00" like the other
# broadcast rows. Fixing this removes the now-unused "00:00:00" shared
# string from the workbook entirely.
#
# Use Formula with a leading apostrophe so the value is stored as literal
# text (matching the existing Text-format / quote-prefixed cell style)
# rather than being reinterpreted as a time value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Formula = "'01:00:00"
$ws.Range("E19").Formula = "'01:00:00"

# Restore the selection that was active when the workbook was last saved.
$ws.Range("E20").Select()
